{"js": "// Rewrite the \"KEY ACHIEVEMENTS AND IMPACT\" bullet list so each bullet is a\n// concise, impact-focused accomplishment statement instead of a job-duty\n// description. Four bullets are rewritten in place and the trailing two\n// bullets (testimony / FEC compliance) are removed entirely.\n//\n// NOTE: several of the original bullet sentences are duplicated verbatim\n// elsewhere in the resume (e.g. inside \"Partner - Siege Analytics\"), so a\n// document-wide text search is not safe. Instead we locate the\n// \"KEY ACHIEVEMENTS AND IMPACT\" heading and only touch the bullet\n// paragraphs that immediately follow it, up to the next section heading.\n\nconst body = context.document.body;\nconst paras = body.paragraphs;\nparas.load(\"items/text,items/style\");\nawait context.sync();\n\n// Locate the \"KEY ACHIEVEMENTS AND IMPACT\" Heading 2 paragraph.\nlet sectionStart = -1;\nfor (let i = 0; i < paras.items.length; i++) {\n  const p = paras.items[i];\n  if (p.style === \"Heading 2\" && p.text.trim() === \"KEY ACHIEVEMENTS AND IMPACT\") {\n    sectionStart = i;\n    break;\n  }\n}\nif (sectionStart === -1) {\n  throw new Error(\"Could not find 'KEY ACHIEVEMENTS AND IMPACT' heading\");\n}\n\n// Find the end of the section: the next Heading 1/Heading 2 paragraph (or\n// end of document).\nlet sectionEnd = paras.items.length;\nfor (let i = sectionStart + 1; i < paras.items.length; i++) {\n  const style = paras.items[i].style;\n  if (style === \"Heading 1\" || style === \"Heading 2\") {\n    sectionEnd = i;\n    break;\n  }\n}\n\n// Old bullet text (trimmed) -> new bullet text, only applied within the\n// KEY ACHIEVEMENTS AND IMPACT section located above.\nconst replacements = {\n  \"\\u2022 Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving classification accuracy from 23% to 64%\":\n    \"\\u2022 Predictive excellence: Achieved 87% voter turnout accuracy vs. 71% industry standard\",\n  \"\\u2022 Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins from \\u00b14.2% to \\u00b12.1%\":\n    \"\\u2022 Reduced polling margins from \\u00b14.2% to \\u00b12.1%\",\n  \"\\u2022 Built redistricting platform used by thousands of analysts nationwide with real-time collaborative editing and Census integration, serving 12,847 analysts across 89 organizations\":\n    \"\\u2022 Executive authority: Briefed Presidents, Congressmen, Senators, Governors on election integrity, voter sentiment and postmortem analysis\",\n  \"\\u2022 Developed longitudinal data analysis methods using geospatial techniques that improved segmentation accuracy by 34% and survey incidence rates by 28%, reducing polling costs while increasing response quality\":\n    \"\\u2022 Methodological advancement: Improved segmentation accuracy 34% and survey incidence 28%\"\n};\n\n// Bullets that should be deleted outright (no replacement).\nconst removals = new Set([\n  \"\\u2022 Provided expert testimony and press briefings on electoral data integrity and demographic modeling accuracy\",\n  \"\\u2022 Demystified FEC compliance through real-time processing systems enabling transparent campaign finance monitoring\"\n]);\n\n// First pass: apply in-place text replacements.\nfor (let i = sectionStart + 1; i < sectionEnd; i++) {\n  const p = paras.items[i];\n  const text = p.text.trim();\n  if (Object.prototype.hasOwnProperty.call(replacements, text)) {\n    p.insertText(replacements[text], Word.InsertLocation.replace);\n  }\n}\nawait context.sync();\n\n// Second pass: delete the two now-obsolete bullet paragraphs. Walk from the\n// end of the section backwards so deleting a paragraph doesn't invalidate\n// the indices of paragraphs we still need to inspect.\nfor (let i = sectionEnd - 1; i > sectionStart; i--) {\n  const p = paras.items[i];\n  if (removals.has(p.text.trim())) {\n    p.delete();\n  }\n}\nawait context.sync();\n", "ps1": "# Rewrite the \"KEY ACHIEVEMENTS AND IMPACT\" bullet list so each bullet is a\n# concise, impact-focused accomplishment statement instead of a job-duty\n# description. Four bullets are rewritten in place and the trailing two\n# bullets (testimony / FEC compliance) are removed entirely.\n#\n# NOTE: several of the original bullet sentences are duplicated verbatim\n# elsewhere in the resume (e.g. inside \"Partner - Siege Analytics\"), so a\n# document-wide Find/Replace is not safe. Instead we locate the\n# \"KEY ACHIEVEMENTS AND IMPACT\" heading paragraph and only touch the bullet\n# paragraphs that immediately follow it, up to the next section heading.\n\n$d = $word.ActiveDocument\n$bullet = [char]0x2022\n$pm = [char]0x00B1\n\n$total = $d.Paragraphs.Count\n\n# Locate the \"KEY ACHIEVEMENTS AND IMPACT\" Heading 2 paragraph.\n$sectionStart = -1\nfor ($i = 1; $i -le $total; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Style.NameLocal -eq \"Heading 2\" -and $p.Range.Text.Trim() -eq \"KEY ACHIEVEMENTS AND IMPACT\") {\n        $sectionStart = $i\n        break\n    }\n}\nif ($sectionStart -eq -1) {\n    throw \"Could not find 'KEY ACHIEVEMENTS AND IMPACT' heading\"\n}\n\n# Find the end of the section: the next Heading 1/Heading 2 paragraph (or\n# end of document).\n$sectionEnd = $total + 1\nfor ($i = $sectionStart + 1; $i -le $total; $i++) {\n    $style = $d.Paragraphs.Item($i).Style.NameLocal\n    if ($style -eq \"Heading 1\" -or $style -eq \"Heading 2\") {\n        $sectionEnd = $i\n        break\n    }\n}\n\n# Old bullet text -> new bullet text (parallel arrays), only applied within\n# the KEY ACHIEVEMENTS AND IMPACT section located above.\n$oldTexts = @(\n    ($bullet + \" Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving classification accuracy from 23% to 64%\"),\n    ($bullet + \" Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins from \" + $pm + \"4.2% to \" + $pm + \"2.1%\"),\n    ($bullet + \" Built redistricting platform used by thousands of analysts nationwide with real-time collaborative editing and Census integration, serving 12,847 analysts across 89 organizations\"),\n    ($bullet + \" Developed longitudinal data analysis methods using geospatial techniques that improved segmentation accuracy by 34% and survey incidence rates by 28%, reducing polling costs while increasing response quality\")\n)\n$newTexts = @(\n    ($bullet + \" Predictive excellence: Achieved 87% voter turnout accuracy vs. 71% industry standard\"),\n    ($bullet + \" Reduced polling margins from \" + $pm + \"4.2% to \" + $pm + \"2.1%\"),\n    ($bullet + \" Executive authority: Briefed Presidents, Congressmen, Senators, Governors on election integrity, voter sentiment and postmortem analysis\"),\n    ($bullet + \" Methodological advancement: Improved segmentation accuracy 34% and survey incidence 28%\")\n)\n\n# Bullets that should be deleted outright (no replacement).\n$removals = @(\n    ($bullet + \" Provided expert testimony and press briefings on electoral data integrity and demographic modeling accuracy\"),\n    ($bullet + \" Demystified FEC compliance through real-time processing systems enabling transparent campaign finance monitoring\")\n)\n\n# First pass: apply in-place text replacements (walk forward, indices are\n# stable because we are only mutating paragraph text, not paragraph count).\nfor ($i = $sectionStart + 1; $i -lt $sectionEnd; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $text = $p.Range.Text.Trim()\n    for ($j = 0; $j -lt $oldTexts.Length; $j++) {\n        if ($text -eq $oldTexts[$j]) {\n            $p.Range.Text = $newTexts[$j]\n        }\n    }\n}\n\n# Second pass: delete the two now-obsolete bullet paragraphs. Walk from the\n# end of the section backwards so deleting a paragraph doesn't invalidate\n# the indices of paragraphs we still need to inspect.\nfor ($i = $sectionEnd - 1; $i -gt $sectionStart; $i--) {\n    $p = $d.Paragraphs.Item($i)\n    $text = $p.Range.Text.Trim()\n    if ($removals -contains $text) {\n        $p.Range.Delete()\n    }\n}\n"}
